$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in D8: "ANTRASİT.jpg" -> "ANRASİT.jpg"
$ws.Range("D8").Value = "ANRASİT.jpg"

# Update selected cell to K19 (cursor position change captured in diff)
$ws.Range("K19").Select()
